$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix exploration time & agent_step_time formula (#13)
# Updated columns: G (Avg_Agent_Step_Time), H (Avg_Experiment_Time),
#                  M (Std_Agent_Step_Time), N (Std_Experiment_Time)
# for data rows 2-13.

$values = @{
    "G2" = 7.40809352;       "H2" = 398.14300002;       "M2" = 1.010627121043755;  "N2" = 105.4411928558437
    "G3" = 11.54658394;      "H3" = 1036.00402306;      "M3" = 1.594986336144876;  "N3" = 293.9491129333546
    "G4" = 2.43512314;       "H4" = 72.60676776;        "M4" = 0.4161116558957283; "N4" = 25.90025464032825
    "G5" = 3.2722633;        "H5" = 163.49937598;       "M5" = 0.562781947459656;  "N5" = 60.38638085824623
    "G6" = 0.7669273600000001; "H6" = 12.3088315;       "M6" = 0.1934415784073938; "N6" = 6.745546862543131
    "G7" = 0.92498132;       "H7" = 24.45588658;        "M7" = 0.171736857617654;  "N7" = 9.945613175554639
    "G8" = 0.40577266;       "H8" = 4.36777126;         "M8" = 0.09406143466749706; "N8" = 2.226559821479375
    "G9" = 0.43350224;       "H9" = 7.997219659999999;  "M9" = 0.09329857909907424; "N9" = 3.844196882024126
    "G10" = 0.2449111;       "H10" = 1.95037218;        "M10" = 0.05633855779818613; "N10" = 0.939673584010758
    "G11" = 0.24452178;      "H11" = 3.54891342;        "M11" = 0.06066215385272521; "N11" = 2.040596071443755
    "G12" = 0.17271368;      "H12" = 1.11245764;        "M12" = 0.04531058831315445; "N12" = 0.6068169002690471
    "G13" = 0.15489318;      "H13" = 1.81991824;        "M13" = 0.0408042252282844;  "N13" = 1.088608759030616
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}
